$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2..83 down to 3..84)
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the latest price data. The leading apostrophe
# forces Excel to store the date as literal text (matching the sheet's
# existing text dates) instead of auto-converting "2026-02-11" into a date
# serial number, while keeping the cell's number format as plain "General"
# like every other cell in the column.
$ws.Cells.Item(2, 1).Value = "'2026-02-11"

$ws.Cells.Item(2, 2).Value = 783.5
$ws.Cells.Item(2, 3).Value = 1112
$ws.Cells.Item(2, 4).Value = 3610
